# Update the "dSF" (column F) values on Sheet1 to reflect the repulled data.
# Mapping of spreadsheet row -> new F value (row number is the Excel row, not the 0-indexed "A" column value).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    2  = 2
    4  = -2
    8  = 1
    15 = 1
    18 = 2
    19 = 1
    31 = 0
    34 = -3
    35 = -3
    37 = 1
    48 = -1
    51 = 1
    60 = -2
    61 = -5
    66 = -2
    68 = -5
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
